$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as serial numbers, matching existing column A values)
$data = @(
    @(44308, 5, 77, 106.3991488068095),
    @(44309, 11, 62, 85.6720418963921),
    @(44310, 14, 62, 85.6720418963921),
    @(44311, 14, 70, 96.7264989152814),
    @(44312, 7, 67, 92.58107753319791)
)

$lastRow = 233
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the formatting of the last existing row's column-A cell (date style)
    # into the new row's column-A cell, then overwrite its value.
    $srcA = $ws.Cells.Item($lastRow, 1)
    $dstA = $ws.Cells.Item($r, 1)
    $srcA.Copy($dstA)
    $dstA.Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
